$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L93").Value = 'fructose, galactose, galactonate, trehalose, starch/glycogen, maltose'
$ws.Range("N93").Value = 'ammonia_assimilation, one nitrite reductase'
$ws.Range("P93").Value = 'sulfate_red_ass'
$ws.Range("Q93").Value = 'branched amino, carbohydrate, LPS export, lipoprotein release, manganese/zinc/iron, molybdate, monosaccharide, oligopeptide, ribose, some type VI secretion'
$ws.Range("L94").Value = 'chitobiose, glycolate'
$ws.Range("Q94").Value = 'heme, iron(III), LPS export, lipoprotein release, phosphate, phospholipid/cholesterol'
$ws.Range("S94").Value = 'one chemotaxis protein'
$ws.Range("T94").Value = 'Oxidative phosphorylation'
$ws.Range("J95").Value = 'formate'
$ws.Range("L95").Value = 'MurNAc, galactose, rhamnulose, cellulose, starch/glycogen, trehalose'
$ws.Range("N95").Value = 'ammonia_assimilation, one nitrite reductase'
$ws.Range("P95").Value = 'sulfate_red_ass, thiosulfate'
$ws.Range("Q95").Value = 'amino acid/amide, branched amino, iron, LPS export, lipoprotein, oligopeptide, phosphate, phospholipid/cholesterol, lipoprotein release, oligopepetide, phosphate, ribose, sodium, urea'
$ws.Range("S95").Value = 'flagellum, chemotaxis (ribose?)'
$ws.Range("T95").Value = 'Oxidative phosphorylation'
$ws.Range("L96").Value = 'chitin, chitobiose, MurNAc, glucose, rhamnulose, galactose, glycolate, glycerate, cellulose, starch/glycogen, maltose'
$ws.Range("N96").Value = 'ammonia_assimilation'
$ws.Range("P96").Value = 'sulfate_red_ass'
$ws.Range("Q96").Value = 'zinc, amino acid/amide, carbohydrate, iron, iron (III), LPS export, lipoprotein release, manganese/iron, molybdate, monosaccharide, oligopeptide, peptide/nickel, phosphate, phospholipid/cholesterol'
$ws.Range("T96").Value = 'Oxidative phosphorylation'
$ws.Range("J97").Value = 'two carotenoid genes'
$ws.Range("L97").Value = 'glucose, alpha-galactosidase, glycolate, maltose, starch/glycogen'
$ws.Range("N97").Value = 'ammonia_assimilation'
$ws.Range("Q97").Value = 'xylose, alpha-glucoside, branched amino, carbohydrate, cobalt/nickel, iron, iron(III), monosaccharide, multiple sugar, phosphate, raffinose/stachyose/melibiose, spermidine/putrescine'
$ws.Range("T97").Value = 'Oxidative phosphorylation'
$ws.Range("J98").Value = 'two carotenoid genes, formate'
$ws.Range("L98").Value = 'chitobiose, glucose, starch/glycogen'
$ws.Range("N98").Value = 'ammonia_assimilation'
$ws.Range("O98").Value = 'partial sulfate_red_ass'
$ws.Range("P98").Value = 'chemotaxis, flagellum'
$ws.Range("Q98").Value = 'cobalt/nickel, iron, iron(III), LPS export, lipoprotein release, microcin C, molybdate, oligopeptide, phosphate, phospholipid/cholesterol, putrescine'
$ws.Range("T98").Value = 'Oxidative phosphorylation'
$ws.Range("J99").Value = 'three carotenoid genes'
$ws.Range("L99").Value = 'chitobiose, cellulose'
$ws.Range("N99").Value = 'ammonia_assimilation, nitronate monooxygenase'
$ws.Range("O99").Value = 'partial sulfate_red_ass'
$ws.Range("Q99").Value = 'heme, LPS export, lipoprotein release, peptide/nickel, phosphate, phospholipid/cholesterol'
$ws.Range("T99").Value = 'Oxidative phosphorylation'
$ws.Range("N100").Value = 'ammonia_assimilation'
$ws.Range("P100").Value = 'sulfate_red_ass'
$ws.Range("Q100").Value = 'iron, iron(III), LPS export, lipoprotein, lipoprotein release, manganese/zinc, monosaccharide, phosphate, phospholipid/cholesterol, ribose, sulfate, urea-binding protein'
$ws.Range("T100").Value = 'Oxidative phosphorylation'
$ws.Range("J101").Value = 'fructose, glucose, rhamnulose, alpha galactosidase, cellulose, cellobiose, trehalose, starch/glycogen, galacturonate'
$ws.Range("P101").Value = 'chemotaxis (monosaccharide, purine, ribose?), some flagellum proteins'
$ws.Range("T101").Value = 'Oxidative phosphorylation'
$ws.Range("J102").Value = 'two carotenoid genes'
$ws.Range("L102").Value = 'chitobiose, glycolate, maltose'
$ws.Range("Q102").Value = 'heme, LPS export, lipoprotein release, phosphate, phospholipid/cholesterol'
$ws.Range("T102").Value = 'Oxidative phosphorylation'
$ws.Range("J103").Value = 'three carotenoid genes'
$ws.Range("L103").Value = 'chitobiose, fructose, rhamnulose, galactose, cellulose degradation, cellobiose, glucoside, starch/glycogen, maltose'
$ws.Range("N103").Value = 'ammonia_assimilation'
$ws.Range("Q103").Value = 'heme, iron(III), LPS export, lipoprotein release, molybdate, phosphate, phospholipid/cholesterol'
$ws.Range("S103").Value = 'one chemotaxis protein'
$ws.Range("T103").Value = 'Oxidative phosphorylation'
$ws.Range("J104").Value = 'methane, methanol, formate, glucose, galacturonate, glucoside, cellobiose, starch/glycogen, maltose, trehalose'
$ws.Range("N104").Value = 'partial nitrate_red_dis, hydroxylamine synthesis, nitrile, ammonia_assimilation'
$ws.Range("P104").Value = 'sulfate_red_ass, alkanesulfonate, methanesulfonate'
$ws.Range("Q104").Value = 'heme, LPS export, lipoprotein release, molybdate, nitrate/nitrite, phosphate, phospholipid/cholesterol, sulfate, sulfonate, urea'
$ws.Range("S104").Value = 'chemotaxis (purine?), flagellum'
$ws.Range("T104").Value = 'Oxidative phosphorylation'
$ws.Range("J105").Value = 'chitobiose, glucose, rhamnulose, glycolate, cellulose, starch/glycogen, maltose, pectin'
$ws.Range("N105").Value = 'ammonia_assimilation'
$ws.Range("P105").Value = 'sulfate_red_ass'
$ws.Range("Q105").Value = 'amino acid/amide, extracellular solute, inositol, iron, iron(III), LPS export, lipoprotein, microcin C, oligopeptide, peptide/nickel, phosphate, phospholipid/cholesterol, xylitol, urea, zinc'
$ws.Range("T105").Value = 'Oxidative phosphorylation'
$ws.Range("N106").Value = 'ammonia_assimilation'
$ws.Range("P106").Value = 'partial sufate_red_ass'
$ws.Range("Q106").Value = 'branched amino, amino acid/amide, biotin, iron, phosphate, polar amino'
$ws.Range("T106").Value = 'Oxidative phosphorylation'
$ws.Range("J107").Value = 'chitobiose, glucose, glucoside, cellobiose'
$ws.Range("K107").Value = 'Wood-Ljungdahl, formate, formaldehyde(?), glycolate'
$ws.Range("N107").Value = 'ammonia_assimilation, nitroalkane, one nitrate_red_ass, one nitrate_red_dis, nitrile'
$ws.Range("P107").Value = 'Sulfur_oxidation(SOX)'
$ws.Range("Q107").Value = 'LPS export, amino acid/amide, branched amino, microcin C, molybdate, nitrate/nitrite, paraquat-inducible protein, phosphate, phospholipid/cholesterol, tungstate'
$ws.Range("T107").Value = 'Oxidative phosphorylation'
$ws.Range("J108").Value = 'formate, glycolate, chitobiose, MurNAc,'
$ws.Range("N108").Value = 'ammonia_assimilation, nitroalkane, transporters for nitrate but no reduction'
$ws.Range("P108").Value = 'sulfate_red_ass, sulfur_oxidation (SOX), taurine, alkanesulfonate, methansulfonate'
$ws.Range("Q108").Value = 'dipthamide, amino acid/amide, branched amino, heme, iron, LPS export, lipoprotein release, molybdate, molbdenum, nitrate/nitrite, phosphate, phospholipid/cholesterol, sulfate, sulfonate, tungstate, urea'
$ws.Range("S108").Value = 'chemotaxis'
$ws.Range("T108").Value = 'Oxidative phosphorylation'
$ws.Range("J109").Value = 'carotenoid synthesis, reductive TCA'
$ws.Range("L109").Value = 'chitobiose, MurNAc, glucose'
$ws.Range("N109").Value = 'nitrogen_fixation, ammonia_assimilation'
$ws.Range("P109").Value = 'sulfate_red_ass, alkanesulfonate'
$ws.Range("Q109").Value = 'amino acid/amide, cobalt/nickel, iron, LPS export, lipopolysaccharide, macrolide, manganese/zinc/iron, molybdate, phosphate, phospholipid/cholesterol, sulfate, sulfonate, thiol reductant'
$ws.Range("S109").Value = 'one chemotaxis protein'
$ws.Range("T109").Value = 'Oxidative phosphorylation'
$ws.Range("J110").Value = 'methanol, formaldehyde, formate, chitobiose, glucose, starch/glycogen, maltose'
$ws.Range("N110").Value = 'ammonia_assimilation, nitrile'
$ws.Range("P110").Value = 'partial sulfate_red_ass, partial sulfur_oxidation (SOX)'
$ws.Range("Q110").Value = 'LPS export, lipoprotein release, molybdate, oligopeptide, phosphate, phospholipid/cholesterol, sulfate'
$ws.Range("T110").Value = 'Oxidative phosphorylation'
$ws.Range("J111").Value = 'chitobiose'
$ws.Range("N111").Value = 'ammonia_assimilation'
$ws.Range("Q111").Value = 'LPS export, lipoprotein release, oligopeptide, phospholipid/cholesterol'
$ws.Range("S111").Value = 'chemotaxis, flagellum'
$ws.Range("T111").Value = 'Oxidative phosphorylation'
$ws.Range("J112").Value = 'carotenoid synthesis, chitobiose, glucose, starch/glycogen, dextrin, maltose'
$ws.Range("L112").Value = 'reductive TCA'
$ws.Range("N112").Value = 'nitrogen_fixation, ammonia_assimilation, hydroxylamine reductase'
$ws.Range("Q112").Value = 'amino acid/amide, capsular polysaccharide, cobalt/nickel, iron, LPS export, lipoprotein release, macrolide, manganese/zinc/iron, molybdate, phosphate, phospholipid/cholesterol, sulfate, sulfonate, type IV secretion'
$ws.Range("T112").Value = 'Oxidative phosphorylation'
$ws.Range("P113").Value = 'sulfate_red_ass, alkanesulfonate'

$ws.Range("J113").Select()
